$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date
$ws.Range("A1").Value = "NBA, Saturday 24th Feb 2024 "

# Update matchup rows with new games/percentages
$ws.Range("A2").Value = "Orlando Magic (31-25) vs Detroit Pistons (8-47)"
$ws.Range("B2").Value = "Orlando Magic (69.23%)"
$ws.Range("C2").Value = "Orlando Magic (70.1%)"

$ws.Range("A3").Value = "Boston Celtics (44-12) vs New York Knicks (34-22)"
$ws.Range("B3").Value = "Boston Celtics (89.66%)"
$ws.Range("C3").Value = "Boston Celtics (72.7%)"

$ws.Range("A4").Value = "Brooklyn Nets (21-34) vs Minnesota Timberwolves (39-17)"
$ws.Range("B4").Value = "Minnesota Timberwolves (76.00%)"
$ws.Range("C4").Value = "Minnesota Timberwolves (79.3%)"

# Remove the now-unused trailing rows (previously rows 5-11)
$ws.Range("A5:C11").Clear()

# Update selected cell to match saved view state
$ws.Range("A4").Select() | Out-Null
